$p = $ppt.ActivePresentation
$lay = $p.SlideMaster.CustomLayouts.Item(4)
for ($j=1; $j -le $lay.Shapes.Count; $j++) {
  $sh = $lay.Shapes.Item($j)
  Write-Output ($j.ToString() + ": " + $sh.Name + " | text=" + $sh.TextFrame.TextRange.Text)
}
